$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Remove the duplicate last row (row 6 was a dupe of row 3: AUS/CDG)
$ws.Rows(6).Delete()

# Refresh the Leave/Return dates from 2013 to 2017 across the remaining data rows
$ws.Range("C2").Value = "12/15/2017"
$ws.Range("D2").Value = "12/31/2017"
$ws.Range("C3").Value = "12/15/2017"
$ws.Range("D3").Value = "12/31/2017"
$ws.Range("C4").Value = "12/15/2017"
$ws.Range("D4").Value = "12/31/2017"
$ws.Range("C5").Value = "12/15/2017"
$ws.Range("D5").Value = "12/31/2017"

# Apply a general number format to the From/To City columns for the data rows
$ws.Range("A2:B5").NumberFormat = "General"

# Move the active selection to D9 (below the now-shorter table)
$ws.Range("D9").Select() | Out-Null
